# Update refreshed M2/FX first- and last-date figures (and a handful of
# M2_Len/FX_Len counters) on the Top33 data-completeness sheet.
# Date columns (E-H) are stored/set as raw Excel serial day numbers since the
# cells already carry the workbook's date number format (style index 2).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 30074
$ws.Range("H2").Value = 45261
$ws.Range("E3").Value = 30011
$ws.Range("F3").Value = 45200
$ws.Range("E4").Value = 30011
$ws.Range("F4").Value = 45200
$ws.Range("G4").Value = 30074
$ws.Range("H4").Value = 45261
$ws.Range("G5").Value = 30074
$ws.Range("H5").Value = 45261
$ws.Range("C6").Value = 443
$ws.Range("F6").Value = 45200
$ws.Range("G6").Value = 30074
$ws.Range("H6").Value = 45261
$ws.Range("G7").Value = 30074
$ws.Range("H7").Value = 45261
$ws.Range("D8").Value = 410
$ws.Range("E8").Value = 30011
$ws.Range("F8").Value = 45200
$ws.Range("H8").Value = 45261
$ws.Range("E9").Value = 30011
$ws.Range("F9").Value = 45200
$ws.Range("G9").Value = 30074
$ws.Range("H9").Value = 45261
$ws.Range("D10").Value = 483
$ws.Range("E10").Value = 30011
$ws.Range("F10").Value = 45200
$ws.Range("H10").Value = 45261
$ws.Range("G11").Value = 30074
$ws.Range("H11").Value = 45261
$ws.Range("C12").Value = 371
$ws.Range("D12").Value = 352
$ws.Range("F12").Value = 45200
$ws.Range("H12").Value = 45261
$ws.Range("C13").Value = 467
$ws.Range("F13").Value = 45200
$ws.Range("G13").Value = 30074
$ws.Range("H13").Value = 45261
$ws.Range("D14").Value = 396
$ws.Range("H14").Value = 45261
$ws.Range("G15").Value = 30042
$ws.Range("H15").Value = 45261
$ws.Range("C16").Value = 455
$ws.Range("D16").Value = 410
$ws.Range("F16").Value = 45200
$ws.Range("H16").Value = 45261
$ws.Range("C17").Value = 370
$ws.Range("D17").Value = 394
$ws.Range("F17").Value = 45200
$ws.Range("H17").Value = 45261
$ws.Range("E18").Value = 30011
$ws.Range("F18").Value = 45200
$ws.Range("G18").Value = 30074
$ws.Range("H18").Value = 45261
$ws.Range("D19").Value = 398
$ws.Range("E19").Value = 30011
$ws.Range("F19").Value = 45200
$ws.Range("H19").Value = 45261
$ws.Range("C20").Value = 481
$ws.Range("F20").Value = 45200
$ws.Range("G20").Value = 30074
$ws.Range("H20").Value = 45261
$ws.Range("C21").Value = 310
$ws.Range("F21").Value = 45200
$ws.Range("G21").Value = 30074
$ws.Range("H21").Value = 45261
$ws.Range("C22").Value = 323
$ws.Range("D22").Value = 367
$ws.Range("F22").Value = 45200
$ws.Range("H22").Value = 45261
$ws.Range("D23").Value = 315
$ws.Range("H23").Value = 45261
$ws.Range("C24").Value = 333
$ws.Range("D24").Value = 313
$ws.Range("F24").Value = 45200
$ws.Range("H24").Value = 45261
$ws.Range("E25").Value = 30011
$ws.Range("F25").Value = 45200
$ws.Range("G25").Value = 30074
$ws.Range("H25").Value = 45261
$ws.Range("D26").Value = 380
$ws.Range("H26").Value = 45261
$ws.Range("D27").Value = 217
$ws.Range("E27").Value = 30011
$ws.Range("F27").Value = 45200
$ws.Range("H27").Value = 45261
$ws.Range("C28").Value = 394
$ws.Range("F28").Value = 45200
$ws.Range("G28").Value = 30074
$ws.Range("H28").Value = 45261
$ws.Range("E29").Value = 30011
$ws.Range("F29").Value = 45200
$ws.Range("G29").Value = 30074
$ws.Range("H29").Value = 45261
$ws.Range("D30").Value = 398
$ws.Range("H30").Value = 45261
$ws.Range("C31").Value = 404
$ws.Range("D31").Value = 317
$ws.Range("F31").Value = 45170
$ws.Range("H31").Value = 45261
$ws.Range("D32").Value = 410
$ws.Range("H32").Value = 45261
$ws.Range("C33").Value = 466
$ws.Range("D33").Value = 317
$ws.Range("F33").Value = 45200
$ws.Range("H33").Value = 45261
$ws.Range("C34").Value = 359
$ws.Range("D34").Value = 367
$ws.Range("F34").Value = 45200
$ws.Range("H34").Value = 45261
